$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.438146710395813
$ws.Range("B1").Value = 3.601300716400146
$ws.Range("C1").Value = 5.548090934753418
$ws.Range("D1").Value = 1.704299569129944
$ws.Range("E1").Value = 0.972909152507782
